$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update status for task 3 (row 6): "In Progress" -> "Completed"
$ws.Range("B6").Value = "Completed"

# Update completion percentage for task 3 (row 6): 0.1 -> 1 (100%)
$ws.Range("D6").Value = 1

# Update the description/note for task 3 (row 6)
$ws.Range("E6").Value = "Includes reading array size metadata and geometry-specific input files. Improvements could be made to pre-processing structure."

# Update the active selection on the sheet
$ws.Range("B7").Select()
